$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.7006386158531
$ws.Range("C2").Value = 8.340294232042817
$ws.Range("D2").Value = 5.9891267849217
$ws.Range("E2").Value = 11.75745654196906
$ws.Range("G2").Value = 30.87871174095526
$ws.Range("H2").Value = 14.8735168670438
$ws.Range("I2").Value = 23.19094960242855
$ws.Range("K2").Value = 8.986505646251979
$ws.Range("L2").Value = 9.702830934424288
$ws.Range("N2").Value = 18.84576132690876
$ws.Range("O2").Value = 22.91385794775655
$ws.Range("B3").Value = 11.43091471675098
$ws.Range("C3").Value = 8.315197352539222
$ws.Range("D3").Value = 5.870759932971582
$ws.Range("E3").Value = 11.76739344145968
$ws.Range("G3").Value = 30.92296222525393
$ws.Range("H3").Value = 14.91498636462818
$ws.Range("I3").Value = 23.26881224354763
$ws.Range("K3").Value = 8.794213910399924
$ws.Range("L3").Value = 9.687258614054372
$ws.Range("N3").Value = 18.90125447841793
$ws.Range("O3").Value = 22.97633073891756
$ws.Range("B4").Value = 11.26397855627403
$ws.Range("C4").Value = 8.299657318748668
$ws.Range("D4").Value = 5.798584840299852
$ws.Range("E4").Value = 11.77559790420314
$ws.Range("G4").Value = 30.95885408212909
$ws.Range("H4").Value = 14.94261315025032
$ws.Range("I4").Value = 23.32064544615681
$ws.Range("K4").Value = 8.675281148054639
$ws.Range("L4").Value = 9.67935985169875
$ws.Range("N4").Value = 18.93694398931173
$ws.Range("O4").Value = 23.01912160866911
$ws.Range("B5").Value = 11.19571755271736
$ws.Range("C5").Value = 8.293293340688134
$ws.Range("D5").Value = 5.769345170151404
$ws.Range("E5").Value = 11.77947054632231
$ws.Range("G5").Value = 30.97566849174305
$ws.Range("H5").Value = 14.95441562731341
$ws.Range("I5").Value = 23.34277950874846
$ws.Range("K5").Value = 8.626664851358299
$ws.Range("L5").Value = 9.676561880412844
$ws.Range("N5").Value = 18.95189539108377
$ws.Range("O5").Value = 23.03767232831805
$ws.Range("B6").Value = 11.18437177540368
$ws.Range("C6").Value = 8.292234761217152
$ws.Range("D6").Value = 5.764501699471739
$ws.Range("E6").Value = 11.78014556998212
$ws.Range("G6").Value = 30.97859250656786
$ws.Range("H6").Value = 14.95640829899304
$ws.Range("I6").Value = 23.34651593520854
$ws.Range("K6").Value = 8.618585086028379
$ws.Range("L6").Value = 9.676122771322586
$ws.Range("N6").Value = 18.95440271423846
$ws.Range("O6").Value = 23.04081984482998
$ws.Range("B7").Value = 11.26305877380682
$ws.Range("C7").Value = 8.299571616607988
$ws.Range("D7").Value = 5.798189745405089
$ws.Range("E7").Value = 11.77564798865148
$ws.Range("G7").Value = 30.95907199462999
$ws.Range("H7").Value = 14.94277011846209
$ws.Range("I7").Value = 23.32093985864259
$ws.Range("K7").Value = 8.674626009516222
$ws.Range("L7").Value = 9.679320409916329
$ws.Range("N7").Value = 18.93714397710656
$ws.Range("O7").Value = 23.01936728529083
$ws.Range("B8").Value = 11.60797147523968
$ws.Range("C8").Value = 8.331668811230553
$ws.Range("D8").Value = 5.948236916065658
$ws.Range("E8").Value = 11.76044655623731
$ws.Range("G8").Value = 30.89215637659805
$ws.Range("H8").Value = 14.88736639969891
$ws.Range("I8").Value = 23.21696087157572
$ws.Range("K8").Value = 8.920422599266638
$ws.Range("L8").Value = 9.697118230639921
$ws.Range("N8").Value = 18.86456056176256
$ws.Range("O8").Value = 22.93447779526555
$ws.Range("B9").Value = 12.26965339889408
$ws.Range("C9").Value = 8.393510830486179
$ws.Range("D9").Value = 6.244488320291051
$ws.Range("E9").Value = 11.74730142001939
$ws.Range("G9").Value = 30.83031891511558
$ws.Range("H9").Value = 14.79588863267688
$ws.Range("I9").Value = 23.04502154973446
$ws.Range("K9").Value = 9.392735119657925
$ws.Range("L9").Value = 9.745080545691508
$ws.Range("N9").Value = 18.73499720461352
$ws.Range("O9").Value = 22.80323985687293
$ws.Range("B10").Value = 12.7415714852155
$ws.Range("C10").Value = 8.438189153087565
$ws.Range("D10").Value = 6.460801144082301
$ws.Range("E10").Value = 11.74776429599195
$ws.Range("G10").Value = 30.82736509242218
$ws.Range("H10").Value = 14.73914027173402
$ws.Range("I10").Value = 22.93821654930768
$ws.Range("K10").Value = 9.730242655210217
$ws.Range("L10").Value = 9.788086931486466
$ws.Range("N10").Value = 18.64751715365927
$ws.Range("O10").Value = 22.72837574286922
$ws.Range("B11").Value = 12.95216331032422
$ws.Range("C11").Value = 8.458331620281868
$ws.Range("D11").Value = 6.558404039975946
$ws.Range("E11").Value = 11.75016076903521
$ws.Range("G11").Value = 30.8352597689276
$ws.Range("H11").Value = 14.71559394584745
$ws.Range("I11").Value = 22.89387350042121
$ws.Range("K11").Value = 9.881025289515051
$ws.Range("L11").Value = 9.809291359814043
$ws.Range("N11").Value = 18.6093781367258
$ws.Range("O11").Value = 22.69901351019997
$ws.Range("B12").Value = 13.03124498814834
$ws.Range("C12").Value = 8.465931249087234
$ws.Range("D12").Value = 6.595208915411497
$ws.Range("E12").Value = 11.75138126735626
$ws.Range("G12").Value = 30.83957686900298
$ws.Range("H12").Value = 14.70700371862657
$ws.Range("I12").Value = 22.87769270826325
$ws.Range("K12").Value = 9.937673932358209
$ws.Range("L12").Value = 9.817552390301119
$ws.Range("N12").Value = 18.59517283211611
$ws.Range("O12").Value = 22.68857074954945
$ws.Range("B13").Value = 13.01424407368958
$ws.Range("C13").Value = 8.464295800872989
$ws.Range("D13").Value = 6.587289848848369
$ws.Range("E13").Value = 11.75110451323655
$ws.Range("G13").Value = 30.83858808863118
$ws.Range("H13").Value = 14.70883926747439
$ws.Range("I13").Value = 22.88115034618544
$ws.Range("K13").Value = 9.925494441856474
$ws.Range("L13").Value = 9.815763014353433
$ws.Range("N13").Value = 18.59822166925405
$ws.Range("O13").Value = 22.6907896999118
$ws.Range("B14").Value = 12.95868315485247
$ws.Range("C14").Value = 8.458957415116373
$ws.Range("D14").Value = 6.561435333179856
$ws.Range("E14").Value = 11.75025491530291
$ws.Range("G14").Value = 30.8355883427932
$ws.Range("H14").Value = 14.71488068404156
$ws.Range("I14").Value = 22.89253004784301
$ws.Range("K14").Value = 9.885695110766934
$ws.Range("L14").Value = 9.809966389664657
$ws.Range("N14").Value = 18.60820471152955
$ws.Range("O14").Value = 22.69814082249752
$ws.Range("B15").Value = 12.92456170592215
$ws.Range("C15").Value = 8.455683816704665
$ws.Range("D15").Value = 6.545577322318992
$ws.Range("E15").Value = 11.7497752321958
$ws.Range("G15").Value = 30.83392375059937
$ws.Range("H15").Value = 14.71862371422054
$ws.Range("I15").Value = 22.89958003211757
$ws.Range("K15").Value = 9.861256802183391
$ws.Range("L15").Value = 9.806445781939022
$ws.Range("N15").Value = 18.61435046191714
$ws.Range("O15").Value = 22.70273167340874
$ws.Range("B16").Value = 12.72771985157871
$ws.Range("C16").Value = 8.436868987366996
$ws.Range("D16").Value = 6.454402989858339
$ws.Range("E16").Value = 11.74765155784395
$ws.Range("G16").Value = 30.82703500734811
$ws.Range("H16").Value = 14.74072473411018
$ws.Range("I16").Value = 22.94119992300836
$ws.Range("K16").Value = 9.720328606759342
$ws.Range("L16").Value = 9.78673380318379
$ws.Range("N16").Value = 18.65004286460148
$ws.Range("O16").Value = 22.73038918408362
$ws.Range("B17").Value = 12.60586036863313
$ws.Range("C17").Value = 8.425278878558958
$ws.Range("D17").Value = 6.398236465955794
$ws.Range("E17").Value = 11.74690773160307
$ws.Range("G17").Value = 30.8251749068171
$ws.Range("H17").Value = 14.75486409288014
$ws.Range("I17").Value = 22.96781971402124
$ws.Range("K17").Value = 9.633129589263586
$ws.Range("L17").Value = 9.775058194686819
$ws.Range("N17").Value = 18.67236242642489
$ws.Range("O17").Value = 22.74855921058457
$ws.Range("B18").Value = 12.53538967142406
$ws.Range("C18").Value = 8.418595534658452
$ws.Range("D18").Value = 6.36585804761576
$ws.Range("E18").Value = 11.74668574139528
$ws.Range("G18").Value = 30.82497472014985
$ws.Range("H18").Value = 14.76321022738555
$ws.Range("I18").Value = 22.98353004546935
$ws.Range("K18").Value = 9.582719049343948
$ws.Range("L18").Value = 9.768497404576886
$ws.Range("N18").Value = 18.68535598868713
$ws.Range("O18").Value = 22.7594518769939
$ws.Range("B19").Value = 12.51146662973397
$ws.Range("C19").Value = 8.416329798498351
$ws.Range("D19").Value = 6.354884013273137
$ws.Range("E19").Value = 11.74664596972531
$ws.Range("G19").Value = 30.82505631847307
$ws.Range("H19").Value = 14.7660727620661
$ws.Range("I19").Value = 22.98891785444026
$ws.Range("K19").Value = 9.565608644265133
$ws.Range("L19").Value = 9.766302739577476
$ws.Range("N19").Value = 18.68978220017681
$ws.Range("O19").Value = 22.76321578102312
$ws.Range("B20").Value = 12.61887246130485
$ws.Range("C20").Value = 8.42651443939355
$ws.Range("D20").Value = 6.404223334546305
$ws.Range("E20").Value = 11.74696561863506
$ws.Range("G20").Value = 30.82528291025134
$ws.Range("H20").Value = 14.75333683236442
$ws.Range("I20").Value = 22.96494465656825
$ws.Range("K20").Value = 9.642438967329953
$ws.Range("L20").Value = 9.776285103189963
$ws.Range("N20").Value = 18.66997033919628
$ws.Range("O20").Value = 22.74657925248322
$ws.Range("B21").Value = 12.97502137542299
$ws.Range("C21").Value = 8.460526199041039
$ws.Range("D21").Value = 6.569033954178124
$ws.Range("E21").Value = 11.7504959793892
$ws.Range("G21").Value = 30.83643342480626
$ws.Range("H21").Value = 14.71309731907429
$ws.Range("I21").Value = 22.88917096604851
$ws.Range("K21").Value = 9.89739774831466
$ws.Range("L21").Value = 9.811662756966607
$ws.Range("N21").Value = 18.60526602299642
$ws.Range("O21").Value = 22.6959632628472
$ws.Range("B22").Value = 13.2038780948888
$ws.Range("C22").Value = 8.482591500719089
$ws.Range("D22").Value = 6.675825367335031
$ws.Range("E22").Value = 11.75462696812131
$ws.Range("G22").Value = 30.85145764147685
$ws.Range("H22").Value = 14.68870013223985
$ws.Range("I22").Value = 22.84320997084465
$ws.Range("K22").Value = 10.0613858361907
$ws.Range("L22").Value = 9.836130428226312
$ws.Range("N22").Value = 18.56435961936377
$ws.Range("O22").Value = 22.66682390925974
$ws.Range("B23").Value = 13.08211450795984
$ws.Range("C23").Value = 8.470830352688562
$ws.Range("D23").Value = 6.61892586933483
$ws.Range("E23").Value = 11.75225579882844
$ws.Range("G23").Value = 30.84273165007086
$ws.Range("H23").Value = 14.70154736855202
$ws.Range("I23").Value = 22.86741408042758
$ws.Range("K23").Value = 9.974120904358633
$ws.Range("L23").Value = 9.822949951788958
$ws.Range("N23").Value = 18.5860660575869
$ws.Range("O23").Value = 22.68201518804273
$ws.Range("B24").Value = 12.61299097173132
$ws.Range("C24").Value = 8.4259559046258
$ws.Range("D24").Value = 6.401516940827739
$ws.Range("E24").Value = 11.74693880725519
$ws.Range("G24").Value = 30.82523137432261
$ws.Range("H24").Value = 14.75402663018635
$ws.Range("I24").Value = 22.96624320463171
$ws.Range("K24").Value = 9.638231061342555
$ws.Range("L24").Value = 9.775729945075129
$ws.Range("N24").Value = 18.67105129790936
$ws.Range("O24").Value = 22.74747300088032
$ws.Range("B25").Value = 12.09278722557736
$ws.Range("C25").Value = 8.37690757309257
$ws.Range("D25").Value = 6.164400315822864
$ws.Range("E25").Value = 11.7490767528817
$ws.Range("G25").Value = 30.83959645641878
$ws.Range("H25").Value = 14.81879836920814
$ws.Range("I25").Value = 22.83496291243497
$ws.Range("K25").Value = 9.266382867891563
$ws.Range("L25").Value = 9.730726065339297
$ws.Range("N25").Value = 18.76868807042954
$ws.Range("O25").Value = 22.83496291243497
